$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second data row (PartCode 124 / Mouse / MTS / ...).
# Use Clear() rather than EntireRow.Delete() -- deleting rows in this
# engine incorrectly decrements any "whole column to the end of the
# sheet" reference (row 1048576 -> 1048575) used by the IsActive data
# validation. Clear() empties the row (which then drops out of
# sheetData/dimension entirely since it becomes fully blank) without
# touching other row references.
$ws.Range("A3:K3").Clear()

# Insert two new columns (ProductMake / BMSMake) between AvailableQty/
# TentativeCost and RGP/IsActive - this shifts the old H:I (RGP,
# IsActive) to J:K.
$ws.Range("H1:I1").EntireColumn.Insert()

# Match the shading/format used by the other PartCode-style columns.
$ws.Range("B1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header row.
$ws.Range("H1").Value2 = "ProductMake"
$ws.Range("I1").Value2 = "BMSMake"

# New data row.
$ws.Range("H2").Value2 = "FBTECH"
$ws.Range("I2").Value2 = "Udaan"

# New column widths.
$ws.Columns.Item(8).ColumnWidth = 11.9
$ws.Columns.Item(9).ColumnWidth = 8.9

# The column insert left the IsActive validation's first area as
# "K1:K3" (row 3 was cleared, not removed, when the column insert ran).
# Dropping validation from the now-empty K3 cell shrinks that area back
# down to K1:K2 without disturbing the J2:J1048576 area.
$ws.Range("K3").Validation.Delete()

$ws.Range("H8").Select()
